# Actualización desde MV -datos-
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct revised figures for row 45 (I-2021) ---
$ws.Range("C45").Value = 299.9
$ws.Range("D45").Value = 1252.6
$ws.Range("E45").Value = 997.8
$ws.Range("F45").Value = 462.5
$ws.Range("G45").Value = 1626.4
$ws.Range("I45").Value = 1007.9
$ws.Range("J45").Value = 634.1

# --- Append new row 48 (III-2021) ---
$ws.Range("A48").Value = "III-2021"
$ws.Range("B48").Value = 8345.200000000001
$ws.Range("C48").Value = 335.7
$ws.Range("D48").Value = 1281.6
$ws.Range("E48").Value = 1035.5
$ws.Range("F48").Value = 478.7
$ws.Range("G48").Value = 1740.4
$ws.Range("H48").Value = 187
$ws.Range("I48").Value = 1100.8
$ws.Range("J48").Value = 649.9
$ws.Range("K48").Value = 1493
$ws.Range("L48").Value = 26.5
$ws.Range("M48").Value = 16.1
